$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Range("P11").Value = 247594
$ws.Range("Q11").Value = 183839
$ws.Range("P13").Value = 1439
$ws.Range("Q13").Value = 1013
$ws.Range("P14").Value = 4814
$ws.Range("Q14").Value = 6287
$ws.Range("P16").Value = 72441
$ws.Range("Q16").Value = 66672
$ws.Range("P18").Value = 2065
$ws.Range("Q18").Value = 568
$ws.Range("F19").Value = 3338
$ws.Range("G19").Value = 2763
$ws.Range("H19").Value = 3644
$ws.Range("I19").Value = 2517
$ws.Range("J19").Value = 3379
$ws.Range("K19").Value = 2190
$ws.Range("L19").Value = 2663
$ws.Range("M19").Value = 3577
$ws.Range("N19").Value = 2488
$ws.Range("O19").Value = 2435
$ws.Range("P19").Value = 2346
$ws.Range("Q19").Value = 1076
$ws.Range("P22").Value = 7942
$ws.Range("Q22").Value = 50668
$ws.Range("F23").Value = 5310
$ws.Range("G23").Value = 4203
$ws.Range("H23").Value = 4433
$ws.Range("I23").Value = 5496
$ws.Range("J23").Value = 4682
$ws.Range("K23").Value = 5305
$ws.Range("L23").Value = 5136
$ws.Range("M23").Value = 7777
$ws.Range("N23").Value = 7930
$ws.Range("O23").Value = 4695
$ws.Range("P23").Value = 4324
$ws.Range("Q23").Value = 5924
$ws.Range("P24").Value = 18
$ws.Range("Q24").Value = 10
$ws.Range("P25").Value = 1390
$ws.Range("Q25").Value = 247
$ws.Range("F26").Value = 13406
$ws.Range("G26").Value = 1286
$ws.Range("H26").Value = 10017
$ws.Range("I26").Value = 5459
$ws.Range("J26").Value = 52003
$ws.Range("K26").Value = 32922
$ws.Range("L26").Value = 28508
$ws.Range("M26").Value = 37909
$ws.Range("N26").Value = 27559
$ws.Range("O26").Value = 37740
$ws.Range("P26").Value = 2256
$ws.Range("Q26").Value = 6237
$ws.Range("P28").Value = 1
$ws.Range("Q28").Value = 1
$ws.Range("P31").Value = 101495
$ws.Range("Q31").Value = 82811
$ws.Range("P32").Value = 66653
$ws.Range("Q32").Value = 57979
$ws.Range("P33").Value = 8318
$ws.Range("Q33").Value = 7868
$ws.Range("P34").Value = 12232
$ws.Range("Q34").Value = 9321
$ws.Range("P35").Value = 30
$ws.Range("Q35").Value = 251
$ws.Range("P36").Value = 44
$ws.Range("Q36").Value = 625
$ws.Range("P37").Value = 362513
$ws.Range("Q37").Value = 485320
$ws.Range("P39").Value = -212973
$ws.Range("Q39").Value = -329501
$ws.Range("P42").Value = 7411
$ws.Range("Q42").Value = 50439
$ws.Range("P43").Value = 16165
$ws.Range("Q43").Value = 58759
$ws.Range("P44").Value = 8754
$ws.Range("Q44").Value = 8320
$ws.Range("P47").Value = 48880
$ws.Range("Q47").Value = 11015
$ws.Range("P49").Value = 18297
$ws.Range("Q49").Value = 0
$ws.Range("L52").Value = 72216
$ws.Range("M52").Value = -73644
$ws.Range("N52").Value = 806
$ws.Range("O52").Value = -104079
$ws.Range("P52").Value = -261000
$ws.Range("Q52").Value = -370042
$ws.Range("F50:K50").Clear()
